$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:G1) is unchanged - "User","Mortal","Fact1".."Fact5".
# Rows 2-7 shift down / get new values (user list fixed: nick/casper/daniel
# swapped out for chelsea/dayna/kelsy, plus an extra "lol" facts row), and
# two new rows (8-9) are appended for personA/personB/personC.
$ws.Cells.Item(2, 1).Value = "personC"
$ws.Cells.Item(2, 2).Value = "praveen"
$ws.Cells.Item(2, 3).Value = "lol"
$ws.Cells.Item(2, 4).Value = "lol2"
$ws.Cells.Item(2, 5).Value = "lol3"
$ws.Cells.Item(2, 6).Value = "lol4"
$ws.Cells.Item(2, 7).Value = "lol5"
$ws.Cells.Item(3, 1).Value = "praveen"
$ws.Cells.Item(3, 2).Value = "joanne"
$ws.Cells.Item(3, 3).Value = "p"
$ws.Cells.Item(3, 4).Value = "2p"
$ws.Cells.Item(3, 5).Value = "3p"
$ws.Cells.Item(3, 6).Value = "4p"
$ws.Cells.Item(3, 7).Value = "5p"
$ws.Cells.Item(4, 1).Value = "joanne"
$ws.Cells.Item(4, 2).Value = "chelsea"
$ws.Cells.Item(4, 3).Value = "j"
$ws.Cells.Item(4, 4).Value = "2j"
$ws.Cells.Item(4, 5).Value = "3j"
$ws.Cells.Item(4, 6).Value = "4j"
$ws.Cells.Item(4, 7).Value = "5j"
$ws.Cells.Item(5, 1).Value = "chelsea"
$ws.Cells.Item(5, 2).Value = "dayna"
$ws.Cells.Item(5, 3).Value = "d"
$ws.Cells.Item(5, 4).Value = "2d"
$ws.Cells.Item(5, 5).Value = "3d"
$ws.Cells.Item(5, 6).Value = "4d"
$ws.Cells.Item(5, 7).Value = "5d"
$ws.Cells.Item(6, 1).Value = "dayna"
$ws.Cells.Item(6, 2).Value = "kelsy"
$ws.Cells.Item(6, 3).Value = "n"
$ws.Cells.Item(6, 4).Value = "2n"
$ws.Cells.Item(6, 5).Value = "3n"
$ws.Cells.Item(6, 6).Value = "4n"
$ws.Cells.Item(6, 7).Value = "5n"
$ws.Cells.Item(7, 1).Value = "kelsy"
$ws.Cells.Item(7, 2).Value = "personA"
$ws.Cells.Item(7, 3).Value = "c"
$ws.Cells.Item(7, 4).Value = "2c"
$ws.Cells.Item(7, 5).Value = "3c"
$ws.Cells.Item(7, 6).Value = "4c"
$ws.Cells.Item(7, 7).Value = "5c"
$ws.Cells.Item(8, 1).Value = "personA"
$ws.Cells.Item(8, 2).Value = "personB"
$ws.Cells.Item(8, 3).Value = "d"
$ws.Cells.Item(8, 4).Value = "2d"
$ws.Cells.Item(8, 5).Value = "3d"
$ws.Cells.Item(8, 6).Value = "4d"
$ws.Cells.Item(8, 7).Value = "5d"
$ws.Cells.Item(9, 1).Value = "personB"
$ws.Cells.Item(9, 2).Value = "personC"
$ws.Cells.Item(9, 3).Value = "d"
$ws.Cells.Item(9, 4).Value = "2d"
$ws.Cells.Item(9, 5).Value = "3d"
$ws.Cells.Item(9, 6).Value = "4d"
$ws.Cells.Item(9, 7).Value = "5d"

# Update the selected cell to match the saved view state
$ws.Range("E13").Select()
